# Scheduled market-data refresh: update computed leve-profit columns (H:N)
# across ALC/ARM/CRP/CUL/LTW/WVR sheets. Values sourced from the latest
# Universalis price snapshot; a handful of rows lost their NQ-profit (M)
# cell entirely because the recipe no longer has an NQ variant, so those
# are cleared instead of overwritten.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 718.2174
$ws.Range("I98").Value = 551.95
$ws.Range("J98").Value = 1826.6666
$ws.Range("K98").Value = 551.95
$ws.Range("L98").Value = 1826.6666
$ws.Range("M98").Value = 946.05
$ws.Range("N98").Value = -4822.6666

# Row 122
$ws.Range("H122").Value = 718.2174
$ws.Range("I122").Value = 551.95
$ws.Range("J122").Value = 1826.6666
$ws.Range("K122").Value = 1655.85
$ws.Range("L122").Value = 5479.9998
$ws.Range("M122").Value = 794.1499999999999
$ws.Range("N122").Value = -10379.9998

# Row 137
$ws.Range("H137").Value = 1980.2632
$ws.Range("I137").Value = 2026.5
$ws.Range("J137").Value = 1786.909
$ws.Range("K137").Value = 6079.5
$ws.Range("L137").Value = 5360.727000000001
$ws.Range("M137").Value = -3529.5
$ws.Range("N137").Value = -10460.727

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2062.54
$ws.Range("I32").Value = 1528.022
$ws.Range("J32").Value = 7467.1113
$ws.Range("K32").Value = 1528.022
$ws.Range("L32").Value = 7467.1113
$ws.Range("M32").Value = -1241.022
$ws.Range("N32").Value = -8041.1113

# Row 74
$ws.Range("H74").Value = 138939.36
$ws.Range("I74").Value = 164645.38
$ws.Range("J74").Value = 56409.527
$ws.Range("K74").Value = 164645.38
$ws.Range("L74").Value = 56409.527
$ws.Range("M74").Value = -163771.38
$ws.Range("N74").Value = -58157.527

# Row 77
$ws.Range("H77").Value = 138939.36
$ws.Range("I77").Value = 164645.38
$ws.Range("J77").Value = 56409.527
$ws.Range("K77").Value = 823226.9
$ws.Range("L77").Value = 282047.635
$ws.Range("M77").Value = -818858.9
$ws.Range("N77").Value = -290783.635

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2391.0557
$ws.Range("I31").Value = 1580.3469
$ws.Range("J31").Value = 4118.2173
$ws.Range("K31").Value = 1580.3469
$ws.Range("L31").Value = 4118.2173
$ws.Range("M31").Value = -1285.3469
$ws.Range("N31").Value = -4708.2173

# Row 34
$ws.Range("H34").Value = 2391.0557
$ws.Range("I34").Value = 1580.3469
$ws.Range("J34").Value = 4118.2173
$ws.Range("K34").Value = 1580.3469
$ws.Range("L34").Value = 4118.2173
$ws.Range("M34").Value = -1378.3469
$ws.Range("N34").Value = -4522.2173

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 934.0345
$ws.Range("I5").Value = 535.1053000000001
$ws.Range("J5").Value = 1692
$ws.Range("K5").Value = 1605.3159
$ws.Range("L5").Value = 5076
$ws.Range("M5").Value = -1493.3159
$ws.Range("N5").Value = -5300

# Row 62
$ws.Range("H62").Value = 2915
$ws.Range("I62").Value = 1495
$ws.Range("J62").Value = 3625
$ws.Range("K62").Value = 4485
$ws.Range("L62").Value = 10875
$ws.Range("M62").Value = -3799
$ws.Range("N62").Value = -12247

# Row 65
$ws.Range("H65").Value = 2915
$ws.Range("I65").Value = 1495
$ws.Range("J65").Value = 3625
$ws.Range("K65").Value = 13455
$ws.Range("L65").Value = 32625
$ws.Range("M65").Value = -10023
$ws.Range("N65").Value = -39489

# Row 70
$ws.Range("H70").Value = 3358.3333
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3358.3333
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10074.9999
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -10704.9999

# Row 73
$ws.Range("H73").Value = 3358.3333
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3358.3333
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10074.9999
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -12258.9999

# Row 74
$ws.Range("H74").Value = 3160.25
$ws.Range("J74").Value = 3160.25
$ws.Range("L74").Value = 9480.75
$ws.Range("N74").Value = -11602.75

# Row 75
$ws.Range("H75").Value = 1576.7333
$ws.Range("I75").Value = 587.6667
$ws.Range("J75").Value = 1824
$ws.Range("K75").Value = 1763.0001
$ws.Range("L75").Value = 5472
$ws.Range("M75").Value = -765.0001
$ws.Range("N75").Value = -7468

# Row 76
$ws.Range("H76").Value = 4166.6665
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4166.6665
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 12499.9995
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -13265.9995

# Row 77
$ws.Range("H77").Value = 3160.25
$ws.Range("J77").Value = 3160.25
$ws.Range("L77").Value = 28442.25
$ws.Range("N77").Value = -39050.25

# Row 78
$ws.Range("H78").Value = 1576.7333
$ws.Range("I78").Value = 587.6667
$ws.Range("J78").Value = 1824
$ws.Range("K78").Value = 5289.0003
$ws.Range("L78").Value = 16416
$ws.Range("M78").Value = -297.0002999999997
$ws.Range("N78").Value = -26400

# Row 79
$ws.Range("H79").Value = 4166.6665
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4166.6665
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 12499.9995
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -15151.9995

# Row 81
$ws.Range("H81").Value = 2848.182
$ws.Range("I81").Value = 333.33334
$ws.Range("J81").Value = 3791.25
$ws.Range("K81").Value = 1000.00002
$ws.Range("L81").Value = 11373.75
$ws.Range("M81").Value = 122.9999799999999
$ws.Range("N81").Value = -13619.75

# Row 82
$ws.Range("H82").Value = 5975.25
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 5975.25
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 17925.75
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -18737.75

# Row 84
$ws.Range("H84").Value = 2848.182
$ws.Range("I84").Value = 333.33334
$ws.Range("J84").Value = 3791.25
$ws.Range("K84").Value = 3000.00006
$ws.Range("L84").Value = 34121.25
$ws.Range("M84").Value = 2615.99994
$ws.Range("N84").Value = -45353.25

# Row 85
$ws.Range("H85").Value = 5975.25
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 5975.25
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 17925.75
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -20733.75

# Row 122
$ws.Range("H122").Value = 1624.875
$ws.Range("I122").Value = 200
$ws.Range("J122").Value = 1828.4286
$ws.Range("K122").Value = 1800
$ws.Range("L122").Value = 16455.8574
$ws.Range("M122").Value = 650
$ws.Range("N122").Value = -21355.8574

# Row 135
$ws.Range("H135").Value = 934.0345
$ws.Range("I135").Value = 535.1053000000001
$ws.Range("J135").Value = 1692
$ws.Range("K135").Value = 4815.947700000001
$ws.Range("L135").Value = 15228
$ws.Range("M135").Value = -2280.947700000001
$ws.Range("N135").Value = -20298

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2544.6594
$ws.Range("I136").Value = 1335.4395
$ws.Range("J136").Value = 5737
$ws.Range("K136").Value = 4006.3185
$ws.Range("L136").Value = 17211
$ws.Range("M136").Value = -1456.3185
$ws.Range("N136").Value = -22311

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1911.1111
$ws.Range("I122").Value = 1671.4286
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 5014.2858
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -2564.2858
$ws.Range("N122").Value = -13150

# Row 126
$ws.Range("H126").Value = 970.2857
$ws.Range("I126").Value = 772.1429000000001
$ws.Range("J126").Value = 1564.7142
$ws.Range("K126").Value = 2316.4287
$ws.Range("L126").Value = 4694.142599999999
$ws.Range("M126").Value = 153.5712999999996
$ws.Range("N126").Value = -9634.142599999999

# Row 132
$ws.Range("H132").Value = 1677.6289
$ws.Range("I132").Value = 1106.16
$ws.Range("J132").Value = 2285.5745
$ws.Range("K132").Value = 3318.48
$ws.Range("L132").Value = 6856.7235
$ws.Range("M132").Value = -788.4800000000005
$ws.Range("N132").Value = -11916.7235

# Row 136
$ws.Range("H136").Value = 10424059
$ws.Range("I136").Value = 13713427
$ws.Range("J136").Value = 418895.62
$ws.Range("K136").Value = 41140281
$ws.Range("L136").Value = 1256686.86
$ws.Range("M136").Value = -41137731
$ws.Range("N136").Value = -1261786.86
